$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Insurance")

# New rows of KPI data to append below the existing Insurance sheet content.
# Columns: A = Category (only on first row of a group), B = KPI name, C = Description
# ht = optional explicit row height (points); when omitted the row keeps the default height.
$rows = @(
  @{ Row=13; A="Claims";    B="Average Cost Per Claim ";            C="In the insurance industry, you are going to have to pay out on claims. That is just the nature of the business. The question is, how much are you going to be paying out? This insurance performance indicator helps estimate this by figuring out the average cost of each claim made"; Ht=30 }
  @{ Row=14;                B="Claim Frequency";                     C=" This key insurance metric measures the likelihood of a loss. It does this by predicting how many claims are to be expected based on the number of policies outstanding. This can help a company manage cashflows, risk exposure, and rate setting."; Ht=30 }
  @{ Row=15;                B="Components of Claim Costs (CCC) ";    C="The CCC metric seeks to provide insight into what costs are associated with a claim. The costs are generally associated with the following items: legal fees, time to settle, administration costs, and report delays. "; Ht=30 }
  @{ Row=16;                B="Average Time to Settle a Claim";      C="he claim settlement time should be used to monitor different policy types as more complicated policies will obviously take longer." }
  @{ Row=17;                B="Client Satisfaction";                 C="Client satisfaction is probably best represented in client retention and policy renewal" }
  @{ Row=18;                B="Calls Handled within 24 Hours";       C="This insurance metric is used to determine how efficient and effective a company’s claims resolution team is. This KPI shouldn’t be used entirely by itself. A company should also consider how many calls the team receives."; Ht=30 }
  @{ Row=19;                B="Underwriting Cycle Time";             C=" This insurance performance indicator measures the number of days it takes the underwriting department of a company to process an insurance policy application. This top insurance KPI can highlight inefficient underwriters, which can have a negative impact on client satisfaction. "; Ht=30 }
  @{ Row=20;                B="Claims Ratio ";                       C="The claims ratio is a very powerful insurance metric. It takes the number of claims made and divides them by the amount of insurance premium earned for a specific period."; Ht=30 }
  @{ Row=21;                C=$null }
  @{ Row=22; A="Financial"; B="Expense Ratio";                       C="The expense ratio performance metric compares the company’s total expenses to the premiums it generates over a specific time period. " }
  @{ Row=23;                B="Loss Ratio";                          C="his insurance KPI divides the total claims payout and divides it by the total premium revenue. A high loss ratio may indicate that policy premiums are set too low."; Ht=30 }
  @{ Row=24;                B="Average Revenue Per Client";          C="We can use this insurance metric to determine the maximum amount of money a company is willing to spend to obtain a new client" }
  @{ Row=25;                B="Cost Per Quote";                      C="The cost per quote takes into consideration all the costs that the company incurs in order to get a quote in front of a potential client." }
  @{ Row=26;                B="Cost Per Bind";                       C="The cost per bind metric determines the incremental cost of binding a new policy. It essentially represents the price a company pays to obtain a new client." }
  @{ Row=27;                B="Net Profit Margin";                   C="This is the official measure of “are you profitable.” If your net income isn’t positive, you aren’t making a profit. But, when you do have a net income that is positive, just divide it by the total revenue."; Ht=30 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    if ($r.ContainsKey("A")) {
        $ws.Cells.Item($rowNum, 1).Value = $r.A
    }
    if ($r.ContainsKey("B")) {
        $ws.Cells.Item($rowNum, 2).Value = $r.B
    }
    if ($r.ContainsKey("C") -and $r.C -ne $null) {
        $c = $ws.Cells.Item($rowNum, 3)
        $c.Value = $r.C
        $c.WrapText = $true
    }
    else {
        $ws.Cells.Item($rowNum, 3).WrapText = $true
    }
    if ($r.ContainsKey("Ht")) {
        $ws.Rows.Item($rowNum).RowHeight = $r.Ht
    }
}

$ws.Range("F16").Select()
$ws.Application.ActiveWindow.ScrollRow = 8
